$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that Excel would otherwise auto-convert
# (numbers, percentages) as literal text, matching the source data which
# stores these columns as plain strings. NumberFormat="@" forces text
# entry; ClearFormats() afterwards drops the temporary formatting so the
# cell style is left exactly as it was (no stray "Text" style/quote-prefix).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") '318.59'
Set-TextValue $ws.Range("E2") '4.26%'

# Row 3
Set-TextValue $ws.Range("D3") '39.83'
Set-TextValue $ws.Range("E3") '2.52%'

# Row 4
Set-TextValue $ws.Range("D4") '5.143'
Set-TextValue $ws.Range("E4") '0.71%'

# Row 5
Set-TextValue $ws.Range("D5") '0.08219'
Set-TextValue $ws.Range("E5") '1.85%'

# Row 6
Set-TextValue $ws.Range("D6") '2.066'
Set-TextValue $ws.Range("E6") '7.08%'

# Row 7
Set-TextValue $ws.Range("D7") '8.386'
Set-TextValue $ws.Range("E7") '4.40%'

# Row 8
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range("D8") '0.9419'
Set-TextValue $ws.Range("E8") '1.69%'

# Row 9
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws.Range("D9") '0.1351'
Set-TextValue $ws.Range("E9") '-5.12%'

# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range("D10") '0.1998'
Set-TextValue $ws.Range("E10") '4.75%'

# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range("D11") '0.09134'
Set-TextValue $ws.Range("E11") '1.21%'

# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range("D12") '0.03524'
Set-TextValue $ws.Range("E12") '0.48%'

# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range("D13") '0.09793'
Set-TextValue $ws.Range("E13") '0.39%'

# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range("D14") '0.001410'
Set-TextValue $ws.Range("E14") '1.04%'

# Row 15
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range("D15") '0.006425'
Set-TextValue $ws.Range("E15") '7.37%'

# Row 16
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range("D16") '3.690'
Set-TextValue $ws.Range("E16") '-1.92%'

# Row 17
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range("D17") '4.319'
Set-TextValue $ws.Range("E17") '2.70%'

# Row 18
Set-TextValue $ws.Range("D18") '3.330'
Set-TextValue $ws.Range("E18") '-1.27%'

# Row 21
Set-TextValue $ws.Range("D21") '4.958'
Set-TextValue $ws.Range("E21") '6.16%'

# Row 22
Set-TextValue $ws.Range("D22") '0.2452'
Set-TextValue $ws.Range("E22") '1.54%'

# Row 23
Set-TextValue $ws.Range("D23") '0.04364'
Set-TextValue $ws.Range("E23") '-0.23%'

# Row 24
Set-TextValue $ws.Range("D24") '0.001236'
Set-TextValue $ws.Range("E24") '2.54%'

# Row 25
Set-TextValue $ws.Range("D25") '0.004789'
Set-TextValue $ws.Range("E25") '12.09%'

# Row 26
Set-TextValue $ws.Range("E26") '-0.02%'

# Row 27
Set-TextValue $ws.Range("D27") '0.0004003'
Set-TextValue $ws.Range("E27") '-10.01%'

# Row 39
Set-TextValue $ws.Range("D39") '0.02349'
Set-TextValue $ws.Range("E39") '15.62%'

# Row 40
Set-TextValue $ws.Range("D40") '0.05210'
Set-TextValue $ws.Range("E40") '3.54%'

# Row 41
Set-TextValue $ws.Range("D41") '0.007758'
Set-TextValue $ws.Range("E41") '3.11%'

# Row 42
Set-TextValue $ws.Range("D42") '0.009881'
Set-TextValue $ws.Range("E42") '1.89%'

# Row 43
Set-TextValue $ws.Range("E43") '4.90%'

# Row 44
Set-TextValue $ws.Range("E44") '-0.57%'

# Row 45
Set-TextValue $ws.Range("D45") '0.009212'
Set-TextValue $ws.Range("E45") '-5.92%'

# Row 46
Set-TextValue $ws.Range("D46") '0.00006610'
Set-TextValue $ws.Range("E46") '6.18%'

# Row 47
Set-TextValue $ws.Range("E47") '-0.06%'

# Row 48
Set-TextValue $ws.Range("D48") '0.002947'
Set-TextValue $ws.Range("E48") '2.53%'

# Row 49
Set-TextValue $ws.Range("E49") '-6.23%'

# Row 50
Set-TextValue $ws.Range("E50") '-0.06%'

# Row 51
Set-TextValue $ws.Range("E51") '-0.06%'
